$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 3: work end time moved from 14:00 to 14:30
$ws.Range("B3").Value = 40592.604166666664

# Row 4: work start time moved from 14:00 to 14:30 (follows row 3's new end time)
$ws.Range("A4").Value = 40592.604166666664

# Row 5: work end time moved from 15:20 to 15:40
$ws.Range("B5").Value = 40593.652777777781

# Row 6: new entry - "Szotar bovitese" task
$ws.Range("A6").Value = 40594.972222222219
$ws.Range("B6").Value = 40595
$ws.Range("C6").Value = "Szótár bővítése"
$ws.Range("D6").Value = "Jégh Tamás, Vad Zsolt"

# Update the active selection shown in the workbook
$ws.Range("G3").Select() | Out-Null
